$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Update financial figures for rows 2-6 (2014-2018 IFRS consolidated) to corrected values

$ws.Range("D2").Value = 7578
$ws.Range("E2").Value = 555
$ws.Range("F2").Value = 599
$ws.Range("G2").Value = 608
$ws.Range("H2").Value = 487
$ws.Range("I2").Value = 419
$ws.Range("J2").Value = 68
$ws.Range("K2").Value = 11109
$ws.Range("L2").Value = 4605
$ws.Range("M2").Value = 6504
$ws.Range("N2").Value = 5479
$ws.Range("O2").Value = 1026
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 709
$ws.Range("R2").Value = -550
$ws.Range("S2").Value = -426
$ws.Range("T2").Value = 308
$ws.Range("U2").Value = 401
$ws.Range("V2").Value = 2178
$ws.Range("W2").Value = 7.32
$ws.Range("X2").Value = 6.43
$ws.Range("Y2").Value = 7.7
$ws.Range("Z2").Value = 4.34
$ws.Range("AA2").Value = 70.79
$ws.Range("AB2").Value = 2569.16
$ws.Range("AC2").Value = 10479
$ws.Range("AD2").Value = 10.97
$ws.Range("AE2").Value = 168577
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 4000
$ws.Range("AH2").Value = 3.48
$ws.Range("AI2").Value = 31.01
$ws.Range("AJ2").Value = 4000000
$ws.Range("D3").Value = 7283
$ws.Range("E3").Value = 259
$ws.Range("F3").Value = 259
$ws.Range("G3").Value = 222
$ws.Range("H3").Value = 171
$ws.Range("I3").Value = 156
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 9037
$ws.Range("L3").Value = 3024
$ws.Range("M3").Value = 6013
$ws.Range("N3").Value = 5128
$ws.Range("O3").Value = 885
$ws.Range("P3").Value = 200
$ws.Range("Q3").Value = 516
$ws.Range("R3").Value = 154
$ws.Range("S3").Value = -496
$ws.Range("T3").Value = 453
$ws.Range("U3").Value = 63
$ws.Range("V3").Value = 1203
$ws.Range("W3").Value = 3.55
$ws.Range("X3").Value = 2.35
$ws.Range("Y3").Value = 2.94
$ws.Range("Z3").Value = 1.7
$ws.Range("AA3").Value = 50.29
$ws.Range("AB3").Value = 2581.07
$ws.Range("AC3").Value = 3904
$ws.Range("AD3").Value = 22.1
$ws.Range("AE3").Value = 157777
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 4000
$ws.Range("AH3").Value = 4.63
$ws.Range("AI3").Value = 83.25
$ws.Range("AJ3").Value = 4000000
$ws.Range("D4").Value = 7145
$ws.Range("E4").Value = 360
$ws.Range("F4").Value = 360
$ws.Range("G4").Value = 427
$ws.Range("H4").Value = 187
$ws.Range("I4").Value = 136
$ws.Range("J4").Value = 52
$ws.Range("K4").Value = 9702
$ws.Range("L4").Value = 3360
$ws.Range("M4").Value = 6342
$ws.Range("N4").Value = 5429
$ws.Range("O4").Value = 913
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 496
$ws.Range("R4").Value = -275
$ws.Range("S4").Value = -139
$ws.Range("T4").Value = 195
$ws.Range("U4").Value = 301
$ws.Range("V4").Value = 1241
$ws.Range("W4").Value = 5.04
$ws.Range("X4").Value = 2.62
$ws.Range("Y4").Value = 2.57
$ws.Range("Z4").Value = 2
$ws.Range("AA4").Value = 52.99
$ws.Range("AB4").Value = 2583.84
$ws.Range("AC4").Value = 3388
$ws.Range("AD4").Value = 25.29
$ws.Range("AE4").Value = 167035
$ws.Range("AF4").Value = 0.51
$ws.Range("AG4").Value = 3500
$ws.Range("AH4").Value = 4.08
$ws.Range("AI4").Value = 83.92
$ws.Range("AJ4").Value = 4000000
$ws.Range("D5").Value = 6943
$ws.Range("E5").Value = -32
$ws.Range("F5").Value = -32
$ws.Range("G5").Value = 319
$ws.Range("H5").Value = 227
$ws.Range("I5").Value = 137
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 9992
$ws.Range("L5").Value = 3456
$ws.Range("M5").Value = 6536
$ws.Range("N5").Value = 5573
$ws.Range("O5").Value = 963
$ws.Range("P5").Value = 200
$ws.Range("Q5").Value = 284
$ws.Range("R5").Value = -378
$ws.Range("S5").Value = 35
$ws.Range("T5").Value = 477
$ws.Range("U5").Value = -193
$ws.Range("V5").Value = 1375
$ws.Range("W5").Value = -0.46
$ws.Range("X5").Value = 3.27
$ws.Range("Y5").Value = 2.5
$ws.Range("Z5").Value = 2.3
$ws.Range("AA5").Value = 52.88
$ws.Range("AB5").Value = 2595.6
$ws.Range("AC5").Value = 3431
$ws.Range("AD5").Value = 25.38
$ws.Range("AE5").Value = 171476
$ws.Range("AF5").Value = 0.51
$ws.Range("AG5").Value = 3500
$ws.Range("AH5").Value = 4.02
$ws.Range("AI5").Value = 82.87
$ws.Range("AJ5").Value = 4000000
$ws.Range("D6").Value = 7892
$ws.Range("E6").Value = 394
$ws.Range("F6").Value = 394
$ws.Range("G6").Value = 712
$ws.Range("H6").Value = 580
$ws.Range("I6").Value = 527
$ws.Range("K6").Value = 10306
$ws.Range("L6").Value = 3369
$ws.Range("M6").Value = 6937
$ws.Range("N6").Value = 5873
$ws.Range("P6").Value = 200
$ws.Range("Q6").Value = 604
$ws.Range("R6").Value = -517
$ws.Range("S6").Value = -56
$ws.Range("T6").Value = 350
$ws.Range("U6").Value = 254
$ws.Range("V6").Value = 1350
$ws.Range("W6").Value = 5
$ws.Range("X6").Value = 7.34
$ws.Range("Y6").Value = 9.21
$ws.Range("Z6").Value = 5.71
$ws.Range("AA6").Value = 48.57
$ws.Range("AB6").Value = 2905.81
$ws.Range("AC6").Value = 13184
$ws.Range("AD6").Value = 6.33
$ws.Range("AE6").Value = 180700
$ws.Range("AF6").Value = 0.46
$ws.Range("AG6").Value = 4000
$ws.Range("AH6").Value = 4.79
$ws.Range("AI6").Value = 24.65
$ws.Range("AJ6").Value = 4000000

# Rows 7-9 (2019E/2020E/2021E forecast rows) had erroneous data - clear all
# financial columns (D:AJ), keeping only A (index), B (연간) and C (period) labels
$ws.Range("D7:AJ9").ClearContents()
